$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay Text,
# matching their original inline-string storage. Without this, values
# such as "0.3800" or "0.000009469" would be auto-coerced by Excel
# into numbers (losing the trailing zero / exact literal form).
# (Each contiguous block is set individually -- this engine only
# honors NumberFormat on the first area of a multi-area Range.)
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7:D16").NumberFormat = "@"
$ws.Range("D18:D27").NumberFormat = "@"
$ws.Range("D29:D43").NumberFormat = "@"
$ws.Range("D45:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.771.30"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "1.886.13"
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "322.54"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "0.4563"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("D8").Value = "0.3800"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").Value = "45.54"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "0.07708"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "0.9606"
$ws.Range("E11").Value = "  -4.05%  "
$ws.Range("D12").Value = "21.97"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "1.882.72"
$ws.Range("E13").Value = "  -5.19%  "
$ws.Range("D14").Value = "6.942"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "5.652"
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("D16").Value = "0.06987"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "83.18"
$ws.Range("E18").Value = "  -6.20%  "
$ws.Range("D19").Value = "0.000009469"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("D20").Value = "16.56"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "28.714.13"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "5.314"
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "2.126.20"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").Value = "2.077"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "155.45"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").Value = "5.591"
$ws.Range("E29").Value = "  -6.83%  "
$ws.Range("D30").Value = "116.82"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "1.799"
$ws.Range("E31").Value = "  -6.16%  "
$ws.Range("D32").Value = "0.09227"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "0.8432"
$ws.Range("E33").Value = "  -5.54%  "
$ws.Range("D34").Value = "5.058"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "1.236"
$ws.Range("E35").Value = "  -8.21%  "
$ws.Range("D36").Value = "2.991"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").Value = "0.05659"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").Value = "1.141"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").Value = "1.001"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "0.02026"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").Value = "7.405"
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("D42").Value = "0.5476"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").Value = "0.1747"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  -29.45%  "
$ws.Range("D45").Value = "9.122"
$ws.Range("E45").Value = "  -7.04%  "
$ws.Range("D46").Value = "2.697"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "11.29"
$ws.Range("E47").Value = "  -6.48%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5143"
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("D49").Value = "0.06797"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("E50").Value = "  -5.57%  "
$ws.Range("D51").Value = "111.43"
$ws.Range("E51").Value = "  -2.54%  "
